$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '92.332.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.113.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '615.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.09'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.397'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.92%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.109.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.734'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.203'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.59%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.232.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.690.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.082.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '449.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000203'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -10.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.264.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("E32").Value = '  -3.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.32'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +22.62%  '
$ws.Range("E35").Value = '  +2.18%  '
$ws.Range("E36").Value = '  -0.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.29'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("E40").Value = '  -3.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '482.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.440'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '160.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("E48").Value = '  -3.64%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0341'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.92%  '
$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.37'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.20'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.31%  '
